$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give E1:L1 the same (date) number format as D1/C1 BEFORE writing the new
# values, so the engine doesn't invent a transient "general date" style for
# them (it only auto-applies one when the existing cell style has no
# explicit number format yet).
$ws.Range("D1").Copy()
$ws.Range("E1:L1").PasteSpecial(-4122) # xlPasteFormats
$ws.Application.CutCopyMode = 0

# Update the date row (row 1, columns D..L)
$ws.Range("D1").Value = (Get-Date -Year 2024 -Month 1 -Day 2 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E1").Value = (Get-Date -Year 2024 -Month 1 -Day 3 -Hour 0 -Minute 0 -Second 0)
$ws.Range("F1").Value = (Get-Date -Year 2024 -Month 1 -Day 4 -Hour 0 -Minute 0 -Second 0)
$ws.Range("G1").Value = (Get-Date -Year 2024 -Month 1 -Day 5 -Hour 0 -Minute 0 -Second 0)
$ws.Range("H1").Value = (Get-Date -Year 2024 -Month 1 -Day 6 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I1").Value = (Get-Date -Year 2024 -Month 1 -Day 7 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J1").Value = (Get-Date -Year 2024 -Month 1 -Day 8 -Hour 0 -Minute 0 -Second 0)
$ws.Range("K1").Value = (Get-Date -Year 2024 -Month 1 -Day 9 -Hour 0 -Minute 0 -Second 0)
$ws.Range("L1").Value = (Get-Date -Year 2024 -Month 1 -Day 10 -Hour 0 -Minute 0 -Second 0)

# Set page orientation to portrait
$ws.PageSetup.Orientation = 1

# Update the selection/active cell
$ws.Range("K10").Select() | Out-Null
